$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) number format on price cells whose new values look numeric,
# so Excel stores them as literal text instead of auto-converting to a Double
# (this matches the original inline-string cell content, e.g. "1.00" must stay "1.00").
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Update the Price (D) and Volume(1h) (E) cell text to the refreshed values.
$ws.Range("D2").Value = '63.622.53'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '3.091.73'
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '592.05'
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").Value = '155.46'
$ws.Range("E6").Value = '  +7.04%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +2.93%  '
$ws.Range("D9").Value = '3.084.82'
$ws.Range("E9").Value = '  -1.57%  '
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = '37.65'
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("E14").Value = '  -1.76%  '
$ws.Range("D15").Value = '3.603.00'
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("D17").Value = '7.20'
$ws.Range("E17").Value = '  -1.51%  '
$ws.Range("D18").Value = '63.589.48'
$ws.Range("D19").Value = '3.089.07'
$ws.Range("E19").Value = '  -1.56%  '
$ws.Range("D20").Value = '476.69'
$ws.Range("E20").Value = '  +1.82%  '
$ws.Range("D21").Value = '14.71'
$ws.Range("E21").Value = '  +2.37%  '
$ws.Range("D22").Value = '0.721'
$ws.Range("E22").Value = '  -1.58%  '
$ws.Range("D23").Value = '7.58'
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("E24").Value = '  +4.36%  '
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("D26").Value = '81.34'
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").Value = '10.03'
$ws.Range("E27").Value = '  +2.84%  '
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '7.40'
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -1.50%  '
$ws.Range("E33").Value = '  +3.81%  '
$ws.Range("D34").Value = '27.33'
$ws.Range("E34").Value = '  -1.84%  '
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("E36").Value = '  -0.86%  '
$ws.Range("E37").Value = '  +6.76%  '
$ws.Range("D38").Value = '6.13'
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("E39").Value = '  -3.01%  '
$ws.Range("D40").Value = '9.35'
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("D41").Value = '50.83'
$ws.Range("E41").Value = '  -1.51%  '
$ws.Range("D42").Value = '444.49'
$ws.Range("E42").Value = '  -1.96%  '
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("D44").Value = '0.0365'
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("E46").Value = '  +3.50%  '
$ws.Range("D47").Value = '2.807.88'
$ws.Range("E47").Value = '  -3.54%  '
$ws.Range("D48").Value = '131.47'
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D49").Value = '25.50'
$ws.Range("E49").Value = '  +5.48%  '
$ws.Range("E51").Value = '  +1.29%  '
